$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = 0.63580578533769971
$ws.Range("AE1").Value = 0.61950929689124923
$ws.Range("AB3").Value = 0.86999522568477672
$ws.Range("AV3").Value = 0.64588504748990672
$ws.Range("L4").Value = 0.95784517020094384
$ws.Range("BJ4").Value = 0.74975495212575904
$ws.Range("BK4").Value = 0.85521885040686851
$ws.Range("S5").Value = 0.74437412247176082
$ws.Range("AE5").Value = 0.99120123522526926
$ws.Range("AM5").Value = 0.98729376573589023
$ws.Range("AX5").Value = 0.93493970772431778
$ws.Range("T6").Value = 0.86160403465120572
$ws.Range("AB6").Value = 0.9631601557256233
$ws.Range("AP6").Value = 0.95526241659065292
$ws.Range("BC7").Value = 0.81388526243621517
$ws.Range("BD7").Value = 0.85265302381906638
$ws.Range("D8").Value = 0.80171983680901149
$ws.Range("G8").Value = 0.70157923852787518
$ws.Range("AR8").Value = 0.64022614397140098
$ws.Range("F9").Value = 0.74541349964024528
$ws.Range("AB10").Value = 0.99794695538076661
$ws.Range("AS10").Value = 0.98535908294806229
$ws.Range("G11").Value = 0.93439944339703551
$ws.Range("L11").Value = 0.69154574791052159
$ws.Range("M11").Value = 0.78499224464105932
$ws.Range("AQ11").Value = 0.59942205363487133
$ws.Range("BM11").Value = 0.99319721178291909
$ws.Range("B12").Value = 0.89033605002675942
$ws.Range("M12").Value = 0.96064885085604035
$ws.Range("O13").Value = 0.86950027973469313
$ws.Range("AG13").Value = 0.70596270127560445
$ws.Range("X14").Value = 0.69696415198683515
$ws.Range("BG15").Value = 0.98666884003661193
$ws.Range("I16").Value = 0.7648983506861915
$ws.Range("AQ16").Value = 0.94399836535183246
$ws.Range("H17").Value = 0.8174090806154275
$ws.Range("R17").Value = 0.6021822251154243
$ws.Range("BM18").Value = 0.83318952423061554
$ws.Range("T19").Value = 0.89039393872770023
$ws.Range("Z19").Value = 0.61945169592079041
$ws.Range("AW19").Value = 0.8659695252219507
$ws.Range("R20").Value = 0.99216642448949965
$ws.Range("AV20").Value = 0.96521051831757032
$ws.Range("F21").Value = 0.75057999664701391
$ws.Range("AS21").Value = 0.73313282544790148
$ws.Range("AL22").Value = 0.70764589759735441
$ws.Range("X23").Value = 0.68400306659200538
$ws.Range("BG23").Value = 0.84376126364458925
$ws.Range("AV24").Value = 0.95114343312243921
$ws.Range("B25").Value = 0.96430097540057458
$ws.Range("BC25").Value = 0.86513138235492293
$ws.Range("BM25").Value = 0.96358358469075445
$ws.Range("G26").Value = 0.57399486759998242
$ws.Range("X26").Value = 0.6226047593517039
$ws.Range("AQ27").Value = 0.83456184379599363
$ws.Range("AY28").Value = 0.98502338476283957
$ws.Range("BE28").Value = 0.824679985930167
$ws.Range("P29").Value = 0.66620399868817826
$ws.Range("S29").Value = 0.8452670291344544
$ws.Range("AF30").Value = 0.96692357488735836
$ws.Range("BH30").Value = 0.9451212088691453
$ws.Range("R31").Value = 0.62219146061300123
$ws.Range("S31").Value = 0.75856126186664941
$ws.Range("AJ32").Value = 0.70756874436351191
$ws.Range("BJ32").Value = 0.94717667151813578
$ws.Range("BO33").Value = 0.8429391237617323
$ws.Range("AF34").Value = 0.60946127060834798
$ws.Range("AJ34").Value = 0.8735210262088331
$ws.Range("AV34").Value = 0.87323315030285298
$ws.Range("U35").Value = 0.95782788645339689
$ws.Range("AR35").Value = 0.80738323390149835
$ws.Range("BN35").Value = 0.92761339793005404
$ws.Range("AL37").Value = 0.80441365220982775
$ws.Range("X39").Value = 0.57786341761320603
$ws.Range("BG39").Value = 0.74711878913027752
$ws.Range("BO39").Value = 0.72033692773356117
$ws.Range("B40").Value = 0.93164200996531576
$ws.Range("I40").Value = 0.77783794354259927
$ws.Range("P41").Value = 0.87367306014909007
$ws.Range("S41").Value = 0.785559848591457
$ws.Range("N42").Value = 0.68217787692398812
$ws.Range("V42").Value = 0.92394843695634898
$ws.Range("Z42").Value = 0.98839998913835692
$ws.Range("AH42").Value = 0.98581504467820791
$ws.Range("AL44").Value = 0.88795028470803539
$ws.Range("BI44").Value = 0.7128556474310811
$ws.Range("S45").Value = 0.73445020783790005
$ws.Range("AU45").Value = 0.62104999996321775
$ws.Range("I46").Value = 0.83749471846330747
$ws.Range("M46").Value = 0.84040251561918322
$ws.Range("U47").Value = 0.97384155835075936
$ws.Range("BP48").Value = 0.87121735400087097
$ws.Range("AD49").Value = 0.96097274428780821
$ws.Range("AU49").Value = 0.83322740729139588
$ws.Range("C50").Value = 0.79613775929948882
$ws.Range("AK50").Value = 0.85816058232773407
$ws.Range("AV50").Value = 0.77303029748548968
$ws.Range("F51").Value = 0.9620747574248929
$ws.Range("N52").Value = 0.84140848201377394
$ws.Range("AV52").Value = 0.74191437197172128
$ws.Range("Z53").Value = 0.94513767404890481
$ws.Range("AN53").Value = 0.82484710894469027
$ws.Range("AY53").Value = 0.96861556618775446
$ws.Range("O54").Value = 0.69674364187853555
$ws.Range("Q54").Value = 0.86373612735124139
$ws.Range("AW54").Value = 0.58883789601451153
$ws.Range("S56").Value = 0.65558029533643036
$ws.Range("X56").Value = 0.76053112864775163
$ws.Range("BC56").Value = 0.96163464868299287
$ws.Range("T57").Value = 0.85113853159728703
$ws.Range("X57").Value = 0.99266859257258777
$ws.Range("AQ58").Value = 0.99592784300150861
$ws.Range("BG58").Value = 0.70417729382194949
$ws.Range("BK58").Value = 0.82240497156025916
$ws.Range("Q59").Value = 0.88523073517414907
$ws.Range("AA59").Value = 0.75226755527300648
$ws.Range("AT59").Value = 0.90016531349516182
$ws.Range("BM59").Value = 0.91346447262331587
$ws.Range("B60").Value = 0.64909502129083918
$ws.Range("T60").Value = 0.91184070414922713
$ws.Range("BJ60").Value = 0.95113906937519954
$ws.Range("BL60").Value = 0.80750005925958535
$ws.Range("M61").Value = 0.86321615753904879
$ws.Range("BK61").Value = 0.55186415196215088
$ws.Range("V62").Value = 0.94975913189967276
$ws.Range("BG64").Value = 0.64585303537549077
$ws.Range("P65").Value = 0.91903701502936253
$ws.Range("AH65").Value = 0.89766692876346921
$ws.Range("AR65").Value = 0.99077642377582475
$ws.Range("BL65").Value = 0.72328467984263467
$ws.Range("AO66").Value = 0.78680116352413032
$ws.Range("P68").Value = 0.87541892670879506
$ws.Range("V68").Value = 0.8442212968056968
$ws.Range("AE68").Value = 0.6158788082800577
$ws.Range("AJ68").Value = 0.99485595967748552
$ws.Range("BO68").Value = 0.73871405636917331
